$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status of second requirement row to "Terminado" (was "Não iniciado")
$ws.Range("E2").Value = "Terminado"

# Move active selection/cursor from B3 to B4
$ws.Range("B4").Select()
